$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: column G "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 7,8,11,12,13,14) {
    $wsOverview.Range("G$r").Value = "2016-08-22 16:22:54"
}

# --- "zh-cn" sheet: column E "Priority" set to "ht", column H "Latest Handoff Datetime" updated ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7,8,11,12,13,14) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-22 16:22:49"
}

# --- "de-de" sheet: column E "Priority" set to "ht", column H "Latest Handoff Datetime" updated ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in 7,8,11,12,13,14) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-22 16:22:54"
}
